$d = $word.ActiveDocument

# Update the date line at the top of the document
$dateRange = $d.Paragraphs(1).Range
$dateRange.Find.Execute("2025-08-21 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-22 Friday", 1) | Out-Null

# Update the division problems inside the table, cell by cell. We scope
# each Find/Replace to the individual cell Range and use wdReplaceOne (1)
# so that cells sharing the same original text (e.g. "80÷8=" appears twice,
# "31÷4=" appears twice) are not all replaced at once.
$t = $d.Tables.Item(1)

$cellRange = $t.Cell(1, 1).Range
$cellRange.Find.Execute("94÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "72÷7=", 1) | Out-Null

$cellRange = $t.Cell(1, 2).Range
$cellRange.Find.Execute("80÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "85÷6=", 1) | Out-Null

$cellRange = $t.Cell(1, 3).Range
$cellRange.Find.Execute("35÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "11÷3=", 1) | Out-Null

$cellRange = $t.Cell(1, 4).Range
$cellRange.Find.Execute("44÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "26÷5=", 1) | Out-Null

$cellRange = $t.Cell(1, 5).Range
$cellRange.Find.Execute("80÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "58÷8=", 1) | Out-Null

$cellRange = $t.Cell(5, 1).Range
$cellRange.Find.Execute("14÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "19÷2=", 1) | Out-Null

$cellRange = $t.Cell(5, 2).Range
$cellRange.Find.Execute("74÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "14÷3=", 1) | Out-Null

$cellRange = $t.Cell(5, 3).Range
$cellRange.Find.Execute("88÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "93÷4=", 1) | Out-Null

$cellRange = $t.Cell(5, 4).Range
$cellRange.Find.Execute("18÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "15÷7=", 1) | Out-Null

$cellRange = $t.Cell(5, 5).Range
$cellRange.Find.Execute("39÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "42÷5=", 1) | Out-Null

$cellRange = $t.Cell(9, 1).Range
$cellRange.Find.Execute("31÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "48÷7=", 1) | Out-Null

$cellRange = $t.Cell(9, 2).Range
$cellRange.Find.Execute("15÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "36÷9=", 1) | Out-Null

$cellRange = $t.Cell(9, 3).Range
$cellRange.Find.Execute("16÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "89÷6=", 1) | Out-Null

$cellRange = $t.Cell(9, 4).Range
$cellRange.Find.Execute("66÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "86÷4=", 1) | Out-Null

$cellRange = $t.Cell(9, 5).Range
$cellRange.Find.Execute("31÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "95÷4=", 1) | Out-Null

$cellRange = $t.Cell(13, 1).Range
$cellRange.Find.Execute("47÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "53÷2=", 1) | Out-Null

$cellRange = $t.Cell(13, 2).Range
$cellRange.Find.Execute("30÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "22÷3=", 1) | Out-Null

$cellRange = $t.Cell(13, 3).Range
$cellRange.Find.Execute("60÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "20÷4=", 1) | Out-Null

$cellRange = $t.Cell(13, 4).Range
$cellRange.Find.Execute("65÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "69÷2=", 1) | Out-Null

$cellRange = $t.Cell(13, 5).Range
$cellRange.Find.Execute("73÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "27÷3=", 1) | Out-Null

$cellRange = $t.Cell(17, 1).Range
$cellRange.Find.Execute("95÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "40÷7=", 1) | Out-Null

$cellRange = $t.Cell(17, 2).Range
$cellRange.Find.Execute("88÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "52÷7=", 1) | Out-Null

$cellRange = $t.Cell(17, 3).Range
$cellRange.Find.Execute("83÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "72÷5=", 1) | Out-Null

$cellRange = $t.Cell(17, 4).Range
$cellRange.Find.Execute("81÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "20÷6=", 1) | Out-Null

$cellRange = $t.Cell(17, 5).Range
$cellRange.Find.Execute("99÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "34÷4=", 1) | Out-Null
